$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Force the literal text into the cell without Excel coercing
    # number-looking strings (e.g. "1.005") into real numbers.
    $escaped = $value.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

$ws.Range("D2").Value = '29.221.38'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '1.826.34'
$ws.Range("E3").Value = '  -0.79%  '

Set-TextCell $ws.Range("D4") '1.005'
$ws.Range("E4").Value = '  +0.55%  '

Set-TextCell $ws.Range("D5") '236.15'
$ws.Range("E5").Value = '  -1.58%  '

Set-TextCell $ws.Range("D6") '0.5948'
$ws.Range("E6").Value = '  -5.04%  '

Set-TextCell $ws.Range("D7") '1.004'
$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws.Range("D8") '0.2810'
$ws.Range("E8").Value = '  -2.75%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws.Range("D9") '0.07019'
$ws.Range("E9").Value = '  -5.19%  '

Set-TextCell $ws.Range("D10") '24.13'
$ws.Range("E10").Value = '  -2.29%  '

Set-TextCell $ws.Range("D11") '0.07671'
$ws.Range("E11").Value = '  -0.78%  '

$ws.Range("D12").Value = '1.829.54'
$ws.Range("E12").Value = '  -0.49%  '

Set-TextCell $ws.Range("D13") '4.730'
$ws.Range("E13").Value = '  -4.90%  '

Set-TextCell $ws.Range("D14") '0.6350'
$ws.Range("E14").Value = '  -6.06%  '

Set-TextCell $ws.Range("D15") '0.000009566'
$ws.Range("E15").Value = '  -5.98%  '

Set-TextCell $ws.Range("D16") '79.14'
$ws.Range("E16").Value = '  -3.40%  '

Set-TextCell $ws.Range("D17") '6.069'
$ws.Range("E17").Value = '  -2.66%  '

$ws.Range("D18").Value = '29.263.19'
$ws.Range("E18").Value = '  -0.11%  '

Set-TextCell $ws.Range("D19") '230.10'
$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range("D20") '1.002'
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws.Range("D21") '11.72'
$ws.Range("E21").Value = '  -4.46%  '

Set-TextCell $ws.Range("D22") '6.943'
$ws.Range("E22").Value = '  -6.27%  '

Set-TextCell $ws.Range("D23") '1.016'
$ws.Range("E23").Value = '  +1.55%  '

Set-TextCell $ws.Range("D24") '156.71'
$ws.Range("E24").Value = '  -1.28%  '

Set-TextCell $ws.Range("D25") '8.040'
$ws.Range("E25").Value = '  -4.89%  '

Set-TextCell $ws.Range("D26") '0.1260'
$ws.Range("E26").Value = '  -6.48%  '

Set-TextCell $ws.Range("D27") '16.47'
$ws.Range("E27").Value = '  -5.20%  '

Set-TextCell $ws.Range("D28") '0.06757'
$ws.Range("E28").Value = '  +1.88%  '

Set-TextCell $ws.Range("D29") '1.464'
$ws.Range("E29").Value = '  +1.13%  '

Set-TextCell $ws.Range("D30") '1.467'
$ws.Range("E30").Value = '  -1.04%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D31") '3.732'
$ws.Range("E31").Value = '  -8.03%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range("D32") '3.679'
$ws.Range("E32").Value = '  -9.38%  '

Set-TextCell $ws.Range("D33") '1.128'
$ws.Range("E33").Value = '  -0.61%  '

Set-TextCell $ws.Range("D34") '1.721'
$ws.Range("E34").Value = '  -5.98%  '

Set-TextCell $ws.Range("D35") '0.6568'
$ws.Range("E35").Value = '  -4.98%  '

Set-TextCell $ws.Range("D36") '2.580'
$ws.Range("E36").Value = '  +0.37%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Range("D37") '2.769'
$ws.Range("E37").Value = '  -1.87%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.224.98'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D39") '0.01759'
$ws.Range("E39").Value = '  -5.15%  '

Set-TextCell $ws.Range("D40") '6.592'
$ws.Range("E40").Value = '  -2.37%  '

Set-TextCell $ws.Range("D41") '0.9340'
$ws.Range("E41").Value = '  -0.03%  '

Set-TextCell $ws.Range("D42") '1.002'
$ws.Range("E42").Value = '  +0.17%  '

$ws.Range("D43").Value = '1.988.99'
$ws.Range("E43").Value = '  +0.87%  '

Set-TextCell $ws.Range("D44") '99.73'
$ws.Range("E44").Value = '  -0.83%  '

Set-TextCell $ws.Range("D45") '63.53'
$ws.Range("E45").Value = '  -2.71%  '

$ws.Range("E46").Value = '  +2.10%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range("D47") '8.738'
$ws.Range("E47").Value = '  -2.90%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range("D48") '1.610'
$ws.Range("E48").Value = '  -5.63%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range("D49") '0.05617'
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws.Range("D50") '0.4577'
$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range("D51") '0.1077'
$ws.Range("E51").Value = '  -6.23%  '
